# The upstream mnova-export code now strips empty NMR datasets before
# writing the workbook, so sheets that only ever contained a header row
# (no data) are removed entirely, and the remaining peak-list sheets are
# regenerated (indices renumbered, one now-empty/duplicate peak row folded
# out of COSY).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the empty datasets (header-only sheets, no peak rows) ---
[void]$wb.Worksheets("H1_1D").Delete()
[void]$wb.Worksheets("H1_pureshift").Delete()
[void]$wb.Worksheets("NOESY").Delete()

# --- COSY: drop a row and fold a duplicate pair, so re-export the table ---
$ws = $wb.Worksheets("COSY")

# Shrink the sheet from 9 rows (header + 8 peaks) to 7 rows (header + 6
# peaks) by deleting the two trailing rows; row count/used-range then
# matches the regenerated export.
$ws.Rows("8:9").Delete()

$cosyData = @(
    @(1, 7.566996222976494,  7.447310850078518,  1.574164390563965,    0),
    @(2, 7.566996222976494,  7.350210639609407,  0.2207812666893005,   0),
    @(3, 3.310011683794853,  2.606743400116869,  0.05913177505135536,  0),
    @(4, 1.965567855644642,  2.605891541783929,  0.3139396905899048,   0),
    @(5, 1.965477164626156,  1.001326068172147,  0.8830317854881287,   0),
    @(6, 1.965097083513452,  1.964488317301542,  0.1004290208220482,   0)
)

for ($i = 0; $i -lt $cosyData.Length; $i++) {
    $row = $i + 2
    for ($c = 0; $c -lt 5; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $cosyData[$i][$c]
    }
}

# Leave the originally-active "molecule" tab selected (sheet deletions /
# edits above shift the active sheet as a side effect; the source diff
# doesn't touch sheet1's view state, so restore it).
$wb.Worksheets("molecule").Activate()
$wb.Worksheets("molecule").Range("A1").Select() | Out-Null
